$wb = $excel.ActiveWorkbook

# Sheet "展览": F3 2786 -> 2806
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F3").Value = 2806

# Sheet "演出": G2 80 -> 180
$wsShow = $wb.Worksheets.Item("演出")
$wsShow.Range("G2").Value = 180

# Sheet "全部类型": G3 80 -> 180, F7 2786 -> 2806
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("G3").Value = 180
$wsAll.Range("F7").Value = 2806
